# Auto-generated from OOXML diff: updates currentAveragePrice-derived market data
# columns: H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
#          K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ
$wb = $excel.ActiveWorkbook

$sheetUpdates = @{
    "ALC" = @(
        @{ Row = 40; Set = @{ "H"=3653.6667; "J"=3653.6667; "L"=3653.6667; "N"=-4003.6667 }; Clear = @() }
        @{ Row = 62; Set = @{ "H"=2409; "J"=2401.5715; "L"=2401.5715; "N"=-3649.5715 }; Clear = @() }
        @{ Row = 64; Set = @{ "H"=0; "I"=0; "J"=0; "K"=0; "L"=0 }; Clear = @("M", "N") }
        @{ Row = 65; Set = @{ "H"=2409; "J"=2401.5715; "L"=12007.8575; "N"=-18247.8575 }; Clear = @() }
        @{ Row = 67; Set = @{ "H"=0; "I"=0; "J"=0; "K"=0; "L"=0 }; Clear = @("M", "N") }
        @{ Row = 95; Set = @{ "H"=50000; "J"=50000; "L"=50000; "N"=-55492 }; Clear = @() }
        @{ Row = 107; Set = @{ "H"=844.5; "I"=894.95; "J"=592.25; "K"=894.95; "L"=592.25; "M"=1025.05; "N"=-4432.25 }; Clear = @() }
        @{ Row = 132; Set = @{ "H"=4540.2144; "I"=3805.3333; "K"=11415.9999; "M"=-8885.999899999999 }; Clear = @() }
        @{ Row = 137; Set = @{ "H"=3339.8286; "I"=2309.5; "K"=6928.5; "M"=-4378.5 }; Clear = @() }
        @{ Row = 138; Set = @{ "H"=1987.123; "I"=1458.8182; "J"=2531.9375; "K"=4376.4546; "L"=7595.8125; "M"=763.5454; "N"=-17875.8125 }; Clear = @() }
    )
    "ARM" = @(
        @{ Row = 74; Set = @{ "H"=2292.3845; "I"=1697; "J"=3632; "K"=1697; "L"=3632; "M"=-823; "N"=-5380 }; Clear = @() }
        @{ Row = 77; Set = @{ "H"=2292.3845; "I"=1697; "J"=3632; "K"=8485; "L"=18160; "M"=-4117; "N"=-26896 }; Clear = @() }
    )
    "BSM" = @(
        @{ Row = 46; Set = @{ "H"=79065; "J"=79065; "L"=79065; "N"=-79661 }; Clear = @() }
        @{ Row = 82; Set = @{ "H"=13200.615; "J"=21665.666; "L"=21665.666; "N"=-22431.666 }; Clear = @() }
        @{ Row = 85; Set = @{ "H"=13200.615; "J"=21665.666; "L"=21665.666; "N"=-24317.666 }; Clear = @() }
        @{ Row = 86; Set = @{ "H"=125002750; "I"=166669400; "J"=2753.5; "K"=166669400; "L"=2753.5; "M"=-166668277; "N"=-4999.5 }; Clear = @() }
        @{ Row = 89; Set = @{ "H"=125002750; "I"=166669400; "J"=2753.5; "K"=833347000; "L"=13767.5; "M"=-833341384; "N"=-24999.5 }; Clear = @() }
        @{ Row = 134; Set = @{ "H"=2385.0952; "I"=1915.9642; "J"=3323.3572; "K"=5747.892599999999; "L"=9970.071599999999; "M"=-3212.892599999999; "N"=-15040.0716 }; Clear = @() }
    )
    "CRP" = @(
        @{ Row = 22; Set = @{ "H"=329.2857; "I"=210; "J"=488.33334; "K"=210; "L"=488.33334; "M"=140; "N"=-1188.33334 }; Clear = @() }
        @{ Row = 31; Set = @{ "H"=4768.4243; "I"=1049.8235; "J"=8719.4375; "K"=1049.8235; "L"=8719.4375; "M"=-754.8235; "N"=-9309.4375 }; Clear = @() }
        @{ Row = 34; Set = @{ "H"=4768.4243; "I"=1049.8235; "J"=8719.4375; "K"=1049.8235; "L"=8719.4375; "M"=-847.8235; "N"=-9123.4375 }; Clear = @() }
        @{ Row = 51; Set = @{ "H"=0; "J"=0; "L"=0 }; Clear = @("N") }
        @{ Row = 61; Set = @{ "H"=0; "J"=0; "L"=0 }; Clear = @("N") }
        @{ Row = 132; Set = @{ "H"=6412349; "I"=1716; "J"=23812638; "K"=5148; "L"=71437914; "M"=-2618; "N"=-71442974 }; Clear = @() }
    )
    "CUL" = @(
        @{ Row = 5; Set = @{ "H"=376.3684; "I"=369.22223; "J"=505; "K"=1107.66669; "L"=1515; "M"=-995.66669; "N"=-1739 }; Clear = @() }
        @{ Row = 75; Set = @{ "H"=3500.5652; "I"=437.66666; "J"=3960; "K"=1312.99998; "L"=11880; "M"=-314.9999800000001; "N"=-13876 }; Clear = @() }
        @{ Row = 78; Set = @{ "H"=3500.5652; "I"=437.66666; "J"=3960; "K"=3938.99994; "L"=35640; "M"=1053.00006; "N"=-45624 }; Clear = @() }
        @{ Row = 80; Set = @{ "H"=3274; "I"=2222; "K"=6666; "M"=-5730 }; Clear = @() }
        @{ Row = 83; Set = @{ "H"=3274; "I"=2222; "K"=19998; "M"=-15318 }; Clear = @() }
        @{ Row = 87; Set = @{ "H"=3302.8; "I"=838; "J"=7000; "K"=2514; "L"=21000; "M"=-1266; "N"=-23496 }; Clear = @() }
        @{ Row = 90; Set = @{ "H"=3302.8; "I"=838; "J"=7000; "K"=7542; "L"=63000; "M"=-1302; "N"=-75480 }; Clear = @() }
        @{ Row = 92; Set = @{ "H"=0; "I"=0; "J"=0; "K"=0; "L"=0 }; Clear = @("M", "N") }
        @{ Row = 135; Set = @{ "H"=376.3684; "I"=369.22223; "J"=505; "K"=3323.00007; "L"=4545; "M"=-788.0000700000001; "N"=-9615 }; Clear = @() }
    )
    "GSM" = @(
        @{ Row = 36; Set = @{ "H"=2500; "I"=0; "J"=2500; "K"=0; "L"=2500; "N"=-3470 }; Clear = @("M") }
    )
    "LTW" = @(
        @{ Row = 136; Set = @{ "H"=20836310; "I"=2560.8; "J"=55559224; "K"=7682.400000000001; "L"=166677672; "M"=-5132.400000000001; "N"=-166682772 }; Clear = @() }
    )
    "WVR" = @(
        @{ Row = 132; Set = @{ "H"=4275792.5; "I"=2387.2778; "J"=7938711.5; "K"=7161.8334; "L"=23816134.5; "M"=-4631.8334; "N"=-23821194.5 }; Clear = @() }
        @{ Row = 136; Set = @{ "H"=2807.5789; "I"=2781.9546; "K"=8345.863799999999; "M"=-5795.863799999999 }; Clear = @() }
    )
}

foreach ($sheetName in $sheetUpdates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($entry in $sheetUpdates[$sheetName]) {
        $rowNum = $entry.Row
        foreach ($col in $entry.Set.Keys) {
            $ws.Range("$col$rowNum").Value = $entry.Set[$col]
        }
        foreach ($col in $entry.Clear) {
            $ws.Range("$col$rowNum").ClearContents()
        }
    }
}